$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New data: the raw GPS export was re-split into training segments ("Trening")
# -- a new column F is added, and the Timestamp column (A) now holds real
# Excel date-time serials (formatted "YYYY-MM-DD HH:MM:SS") instead of plain
# text strings. Rows 8-13 are brand-new data points for the "Mała Gra" part
# of the session.
# ---------------------------------------------------------------------------

$rows = @(
  @{ Row=2;  A=45686.47739340278; B=1294.5; C=10.19; D=0.9392164945602425; E="10-15"; F="Duża Gra" }
  @{ Row=3;  A=45686.47775335648; B=1325.6; C=10.18; D=0.8800620010920935; E="10-15"; F="Duża Gra" }
  @{ Row=4;  A=45686.47835868056; B=1377.9; C=10.33; D=1.952321222850257;  E="10-15"; F="Duża Gra" }
  @{ Row=5;  A=45686.47475682871; B=1066.7; C=5.19;  D=1.612820591245379;  E="5-10";  F="Duża Gra" }
  @{ Row=6;  A=45686.47833206019; B=1375.6; C=8.08;  D=1.756553990500314;  E="5-10";  F="Duża Gra" }
  @{ Row=7;  A=45686.47835752315; B=1377.8; C=9.800000000000001; D=2.019265532493593; E="5-10"; F="Duża Gra" }
  @{ Row=8;  A=45686.48793854166; B=2205.6; C=12.68; D=2.935751097542898; E="10-15"; F="Mała Gra" }
  @{ Row=9;  A=45686.4937869213;  B=2710.9; C=12.43; D=3.117605243410381; E="10-15"; F="Mała Gra" }
  @{ Row=10; A=45686.49390960648; B=2721.5; C=12.78; D=3.317203487668718; E="10-15"; F="Mała Gra" }
  @{ Row=11; A=45686.48218506944; B=1708.5; C=8.93;  D=2.816275221960884; E="5-10";  F="Mała Gra" }
  @{ Row=12; A=45686.49378344908; B=2710.6; C=9.19;  D=2.775653447423661; E="5-10";  F="Mała Gra" }
  @{ Row=13; A=45686.49546747685; B=2856.1; C=9.81;  D=2.644420794078283; E="5-10";  F="Mała Gra" }
)

# New header for column F, styled like the rest of the header row (bold,
# bordered, centered) by copying the format from E1.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Trening"

# Give column A a real datetime number format. Apply it to the first data
# cell in two steps (lower-case, then upper-case form) so the workbook keeps
# both format codes registered exactly like the source file, then reuse the
# resulting style for the remaining rows.
$firstCell = $ws.Cells.Item(2, 1)
$firstCell.NumberFormat = "yyyy-mm-dd h:mm:ss"
$firstCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

foreach ($r in $rows) {
    if ($r.Row -gt 2) {
        $ws.Cells.Item($r.Row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
